# "Generate Report for handoff"
#
# The handoff-status workbook is regenerated for a new source file
# (28a65812-a08d-4615-9da0-bd8d5a15fcfc -> 98539602-0ccb-4ca6-9ca8-8ca764535844,
# with a new handoff-package hash cb5e1bf0... -> 732445440...) and the stale
# "47dbaa12.../Handoff transform failed" row is dropped from every sheet
# now that handoff succeeded.

$wb = $excel.ActiveWorkbook

$oldGuid = "28a65812-a08d-4615-9da0-bd8d5a15fcfc"
$newGuid = "98539602-0ccb-4ca6-9ca8-8ca764535844"
$oldHash = "cb5e1bf09579b5fffbc4bfccf62dfdbf137dfa97"
$newHash = "732445440fb9500936fd99527b2b9cc8cc6a20f5"

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# ---- 1. Update the source-file identifiers / hashes / timestamps in place ----

$ws1.Range("A2").Value = $newGuid + ".md"

$ws2.Range("A2").Value = $newGuid + ".md"
$ws2.Range("C2").Value = $newGuid + "." + $newHash + ".zh-cn.xlf"
$ws2.Range("D2").Value = "2016-01-18 04:02:11"

$ws3.Range("A2").Value = $newGuid + ".md"
$ws3.Range("C2").Value = $newGuid + "." + $newHash + ".de-de.xlf"
$ws3.Range("D2").Value = "2016-01-18 04:02:24"

# ---- 2. Drop the "Handoff transform failed" row (row 3) from every sheet ----
# (the remaining ".localization-config" row shifts up to take its place)

$ws1.Rows.Item(3).Delete()
$ws2.Rows.Item(3).Delete()
$ws3.Rows.Item(3).Delete()

# ---- 3. Rebuild the hyperlinks so their display text / targets follow the
#         renamed file and the row shift above ----

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/ea3d4cf4c78b0befa7609568df1ce3a7be479ec3/e2e/" + $newGuid + ".md", "", "", $newGuid + ".md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ea3d4cf4c78b0befa7609568df1ce3a7be479ec3/.localization-config", "", "", ".localization-config") | Out-Null

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/ea3d4cf4c78b0befa7609568df1ce3a7be479ec3/e2e/" + $newGuid + ".md", "", "", $newGuid + ".md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e661027ea602de4c59eccdd04eef4670491378bc/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/" + $newGuid + "." + $newHash + ".zh-cn.xlf", "", "", $newGuid + "." + $newHash + ".zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ea3d4cf4c78b0befa7609568df1ce3a7be479ec3/.localization-config", "", "", ".localization-config") | Out-Null

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/ea3d4cf4c78b0befa7609568df1ce3a7be479ec3/e2e/" + $newGuid + ".md", "", "", $newGuid + ".md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/50f0aa75b7a82cb8047177d40dad97be1c0d5a2e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/" + $newGuid + "." + $newHash + ".de-de.xlf", "", "", $newGuid + "." + $newHash + ".de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ea3d4cf4c78b0befa7609568df1ce3a7be479ec3/.localization-config", "", "", ".localization-config") | Out-Null

"Report regenerated for " + $newGuid
